$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-66) holds an "updated" date stamp that is bumped
# by one day (45179 -> 45180) on every automatic refresh.
$ws.Range("C2:C66").Value = 45180
